$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column J, matching the style/formatting of the existing
# header row (copy I1's format onto J1, then set the text).
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "Request Frequency"

# Row 2 stays blank/empty like the rest of that row, but the cell itself
# still needs to exist in the sheet (matching A2:I2). Touch a formatting
# property with its own default value so a cell record is materialized
# without leaving any value or picking up a new/non-default style.
$ws.Range("J2").Font.Bold = $false

# New data values for the existing config rows
$ws.Range("J3").Value = 0.9
$ws.Range("J4").Value = 0.2
